$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "BTC"
$ws.Range("C2").Value = "Bitcoin"
$ws.Range("D2").Value = 26504
$ws.Range("E2").Value = 514130026129
$ws.Range("F2").Value = 7719628543
$ws.Range("G2").Value = -0.00517

$ws.Range("B3").Value = "ETH"
$ws.Range("C3").Value = "Ethereum"
$ws.Range("D3").Value = 1733.98
$ws.Range("E3").Value = 208408939062
$ws.Range("F3").Value = 3500759457
$ws.Range("G3").Value = -0.20914

$ws.Range("B4").Value = "USDT"
$ws.Range("C4").Value = "Tether"
$ws.Range("D4").Value = 0.999698
$ws.Range("E4").Value = 83123063318
$ws.Range("F4").Value = 11356066225
$ws.Range("G4").Value = 0.0283

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "BNB"
$ws.Range("D5").Value = 246.97
$ws.Range("E5").Value = 38451250501
$ws.Range("F5").Value = 516731160
$ws.Range("G5").Value = 0.64275

$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "USD Coin"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 28343670246
$ws.Range("F6").Value = 1776435711
$ws.Range("G6").Value = -0.01725

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "XRP"
$ws.Range("D7").Value = 0.489456
$ws.Range("E7").Value = 25416120120
$ws.Range("F7").Value = 576976568
$ws.Range("G7").Value = 1.84822

$ws.Range("B8").Value = "STETH"
$ws.Range("C8").Value = "Lido Staked Ether"
$ws.Range("D8").Value = 1734.03
$ws.Range("E8").Value = 12560665943
$ws.Range("F8").Value = 10084881
$ws.Range("G8").Value = -0.17523

$ws.Range("B9").Value = "ADA"
$ws.Range("C9").Value = "Cardano"
$ws.Range("D9").Value = 0.266724
$ws.Range("E9").Value = 9339948530
$ws.Range("F9").Value = 138467355
$ws.Range("G9").Value = -0.57603

$ws.Range("B10").Value = "DOGE"
$ws.Range("C10").Value = "Dogecoin"
$ws.Range("D10").Value = 0.062283
$ws.Range("E10").Value = 8705156545
$ws.Range("F10").Value = 187490734
$ws.Range("G10").Value = -0.2695

$ws.Range("B11").Value = "TRX"
$ws.Range("C11").Value = "TRON"
$ws.Range("D11").Value = 0.070606
$ws.Range("E11").Value = 6360975725
$ws.Range("F11").Value = 241906419
$ws.Range("G11").Value = -0.8918199999999999

$ws.Range("B12").Value = "SOL"
$ws.Range("C12").Value = "Solana"
$ws.Range("D12").Value = 15.73
$ws.Range("E12").Value = 6288152256
$ws.Range("F12").Value = 148119083
$ws.Range("G12").Value = -0.44167

$ws.Range("B13").Value = "DOT"
$ws.Range("C13").Value = "Polkadot"
$ws.Range("D13").Value = 4.61
$ws.Range("E13").Value = 5727178003
$ws.Range("F13").Value = 90357091
$ws.Range("G13").Value = 1.30981

$ws.Range("B14").Value = "MATIC"
$ws.Range("C14").Value = "Polygon"
$ws.Range("D14").Value = 0.610446
$ws.Range("E14").Value = 5667152605
$ws.Range("F14").Value = 234043128
$ws.Range("G14").Value = -1.37483

$ws.Range("B15").Value = "LTC"
$ws.Range("C15").Value = "Litecoin"
$ws.Range("D15").Value = 77.36
$ws.Range("E15").Value = 5661924459
$ws.Range("F15").Value = 226898136
$ws.Range("G15").Value = 0.39579

$ws.Range("B16").Value = "DAI"
$ws.Range("C16").Value = "Dai"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 4416010684
$ws.Range("F16").Value = 78231889
$ws.Range("G16").Value = -0.08398

$ws.Range("B17").Value = "BUSD"
$ws.Range("C17").Value = "Binance USD"
$ws.Range("D17").Value = 0.999952
$ws.Range("E17").Value = 4293761136
$ws.Range("F17").Value = 1165530473
$ws.Range("G17").Value = -0.0453

$ws.Range("B18").Value = "SHIB"
$ws.Range("C18").Value = "Shiba Inu"
$ws.Range("D18").Value = 0.00000724
$ws.Range("E18").Value = 4254974703
$ws.Range("F18").Value = 129185296
$ws.Range("G18").Value = 4.99518

$ws.Range("B19").Value = "WBTC"
$ws.Range("C19").Value = "Wrapped Bitcoin"
$ws.Range("D19").Value = 26508
$ws.Range("E19").Value = 4157107522
$ws.Range("F19").Value = 57980137
$ws.Range("G19").Value = -0.05611

$ws.Range("B20").Value = "AVAX"
$ws.Range("C20").Value = "Avalanche"
$ws.Range("D20").Value = 11.52
$ws.Range("E20").Value = 3975061234
$ws.Range("F20").Value = 102623100
$ws.Range("G20").Value = -2.10462

$ws.Range("B21").Value = "UNI"
$ws.Range("C21").Value = "Uniswap"
$ws.Range("D21").Value = 4.54
$ws.Range("E21").Value = 3418290137
$ws.Range("F21").Value = 28890358
$ws.Range("G21").Value = -1.2367

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "LEO Token"
$ws.Range("D22").Value = 3.55
$ws.Range("E22").Value = 3301646747
$ws.Range("F22").Value = 485957
$ws.Range("G22").Value = 0.94928

$ws.Range("B23").Value = "TUSD"
$ws.Range("C23").Value = "TrueUSD"
$ws.Range("D23").Value = 1.002
$ws.Range("E23").Value = 3116271882
$ws.Range("F23").Value = 112738765
$ws.Range("G23").Value = 0.02121

$ws.Range("B24").Value = "LINK"
$ws.Range("C24").Value = "Chainlink"
$ws.Range("D24").Value = 5.24
$ws.Range("E24").Value = 2710821626
$ws.Range("F24").Value = 80609426
$ws.Range("G24").Value = -2.04024

$ws.Range("B25").Value = "ATOM"
$ws.Range("C25").Value = "Cosmos Hub"
$ws.Range("D25").Value = 8.77
$ws.Range("E25").Value = 2559075861
$ws.Range("F25").Value = 56145492
$ws.Range("G25").Value = -1.26254

$ws.Range("B26").Value = "XMR"
$ws.Range("C26").Value = "Monero"
$ws.Range("D26").Value = 139.38
$ws.Range("E26").Value = 2532899268
$ws.Range("F26").Value = 46699578
$ws.Range("G26").Value = 2.77859

$ws.Range("B27").Value = "OKB"
$ws.Range("C27").Value = "OKB"
$ws.Range("D27").Value = 41.5
$ws.Range("E27").Value = 2492565490
$ws.Range("F27").Value = 2660769
$ws.Range("G27").Value = -0.73752

$ws.Range("B28").Value = "ETC"
$ws.Range("C28").Value = "Ethereum Classic"
$ws.Range("D28").Value = 15.45
$ws.Range("E28").Value = 2184134520
$ws.Range("F28").Value = 60144446
$ws.Range("G28").Value = 0.65201

$ws.Range("B29").Value = "XLM"
$ws.Range("C29").Value = "Stellar"
$ws.Range("D29").Value = 0.080692
$ws.Range("E29").Value = 2171511121
$ws.Range("F29").Value = 32287669
$ws.Range("G29").Value = 2.60238

$ws.Range("B30").Value = "BCH"
$ws.Range("C30").Value = "Bitcoin Cash"
$ws.Range("D30").Value = 108.06
$ws.Range("E30").Value = 2100080215
$ws.Range("F30").Value = 48058055
$ws.Range("G30").Value = 1.05591

$ws.Range("B31").Value = "TON"
$ws.Range("C31").Value = "Toncoin"
$ws.Range("D31").Value = 1.41
$ws.Range("E31").Value = 2081536704
$ws.Range("F31").Value = 3149825
$ws.Range("G31").Value = -0.53277

$ws.Range("B32").Value = "ICP"
$ws.Range("C32").Value = "Internet Computer"
$ws.Range("D32").Value = 3.99
$ws.Range("E32").Value = 1740630565
$ws.Range("F32").Value = 9826592
$ws.Range("G32").Value = -0.29666

$ws.Range("B33").Value = "FIL"
$ws.Range("C33").Value = "Filecoin"
$ws.Range("D33").Value = 3.69
$ws.Range("E33").Value = 1589042737
$ws.Range("F33").Value = 54585734
$ws.Range("G33").Value = -1.05203

$ws.Range("B34").Value = "LDO"
$ws.Range("C34").Value = "Lido DAO"
$ws.Range("D34").Value = 1.78
$ws.Range("E34").Value = 1566185230
$ws.Range("F34").Value = 15396497
$ws.Range("G34").Value = -1.17958

$ws.Range("B35").Value = "QNT"
$ws.Range("C35").Value = "Quant"
$ws.Range("D35").Value = 101.39
$ws.Range("E35").Value = 1468856926
$ws.Range("F35").Value = 46266445
$ws.Range("G35").Value = -10.24997

$ws.Range("B36").Value = "HBAR"
$ws.Range("C36").Value = "Hedera"
$ws.Range("D36").Value = 0.04587823
$ws.Range("E36").Value = 1447494423
$ws.Range("F36").Value = 14169289
$ws.Range("G36").Value = 0.3569

$ws.Range("B37").Value = "APT"
$ws.Range("C37").Value = "Aptos"
$ws.Range("D37").Value = 6.98
$ws.Range("E37").Value = 1438456139
$ws.Range("F37").Value = 134118357
$ws.Range("G37").Value = 1.70374

$ws.Range("B38").Value = "CRO"
$ws.Range("C38").Value = "Cronos"
$ws.Range("D38").Value = 0.053854
$ws.Range("E38").Value = 1406369027
$ws.Range("F38").Value = 4073055
$ws.Range("G38").Value = 0.9067

$ws.Range("B39").Value = "ARB"
$ws.Range("C39").Value = "Arbitrum"
$ws.Range("D39").Value = 1.013
$ws.Range("E39").Value = 1292434959
$ws.Range("F39").Value = 97846537
$ws.Range("G39").Value = 1.34716

$ws.Range("B40").Value = "NEAR"
$ws.Range("C40").Value = "NEAR Protocol"
$ws.Range("D40").Value = 1.25
$ws.Range("E40").Value = 1152941419
$ws.Range("F40").Value = 34779761
$ws.Range("G40").Value = -0.42269

$ws.Range("B41").Value = "VET"
$ws.Range("C41").Value = "VeChain"
$ws.Range("D41").Value = 0.01510038
$ws.Range("E41").Value = 1097545815
$ws.Range("F41").Value = 32461668
$ws.Range("G41").Value = 0.11249

$ws.Range("B42").Value = "USDP"
$ws.Range("C42").Value = "Pax Dollar"
$ws.Range("D42").Value = 1.001
$ws.Range("E42").Value = 1007076455
$ws.Range("F42").Value = 1216097
$ws.Range("G42").Value = 0.01424

$ws.Range("B43").Value = "FRAX"
$ws.Range("C43").Value = "Frax"
$ws.Range("D43").Value = 0.9991370000000001
$ws.Range("E43").Value = 1002783822
$ws.Range("F43").Value = 13972447
$ws.Range("G43").Value = -0.10188

$ws.Range("B44").Value = "GRT"
$ws.Range("C44").Value = "The Graph"
$ws.Range("D44").Value = 0.101666
$ws.Range("E44").Value = 915694293
$ws.Range("F44").Value = 35304813
$ws.Range("G44").Value = 1.51025

$ws.Range("B45").Value = "ALGO"
$ws.Range("C45").Value = "Algorand"
$ws.Range("D45").Value = 0.118602
$ws.Range("E45").Value = 858929426
$ws.Range("F45").Value = 24930131
$ws.Range("G45").Value = -1.56302

$ws.Range("B46").Value = "RETH"
$ws.Range("C46").Value = "Rocket Pool ETH"
$ws.Range("D46").Value = 1864.18
$ws.Range("E46").Value = 833174097
$ws.Range("F46").Value = 2068802
$ws.Range("G46").Value = -0.15891

$ws.Range("B47").Value = "RPL"
$ws.Range("C47").Value = "Rocket Pool"
$ws.Range("D47").Value = 41.21
$ws.Range("E47").Value = 801936794
$ws.Range("F47").Value = 3947950
$ws.Range("G47").Value = -1.19626

$ws.Range("B48").Value = "STX"
$ws.Range("C48").Value = "Stacks"
$ws.Range("D48").Value = 0.572309
$ws.Range("E48").Value = 793659211
$ws.Range("F48").Value = 26127467
$ws.Range("G48").Value = -3.48531

$ws.Range("B49").Value = "APE"
$ws.Range("C49").Value = "ApeCoin"
$ws.Range("D49").Value = 2.13
$ws.Range("E49").Value = 784794655
$ws.Range("F49").Value = 58858802
$ws.Range("G49").Value = 0.37583

$ws.Range("B50").Value = "EGLD"
$ws.Range("C50").Value = "MultiversX"
$ws.Range("D50").Value = 30.59
$ws.Range("E50").Value = 782089156
$ws.Range("F50").Value = 9191552
$ws.Range("G50").Value = -0.27369

$ws.Range("B51").Value = "FTM"
$ws.Range("C51").Value = "Fantom"
$ws.Range("D51").Value = 0.273029
$ws.Range("E51").Value = 762847483
$ws.Range("F51").Value = 76183711
$ws.Range("G51").Value = 2.60238
